$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels from uppercase to title case
$ws.Range("A1").Value = "Question"
$ws.Range("B1").Value = "Answer"

# Update the active cell selection
$ws.Range("D11").Select()
